$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 61, shifting existing rows 61:86 down to 62:87.
$ws.Rows(61).Insert()

# Populate the new row 61 with the new weekly record (same fixed attributes
# as the surrounding Tuna / Vega Modelo de Temuco rows, new date + price data).
$ws.Range("A61").Value = 10
$ws.Range("B61").Value = "Vega Modelo de Temuco"
$ws.Range("C61").Value = "La Araucanía"
$ws.Range("D61").Value = 45029
$ws.Range("E61").Value = 9
$ws.Range("F61").Value = "Fruta"
$ws.Range("G61").Value = 100107
$ws.Range("H61").Value = "Otros"
$ws.Range("I61").Value = 100107011
$ws.Range("J61").Value = "Tuna"
$ws.Range("K61").Value = "Sin especificar"
$ws.Range("L61").Value = "Primera"
$ws.Range("M61").Value = 35
$ws.Range("N61").Value = 22000
$ws.Range("O61").Value = 22000
$ws.Range("P61").Value = 22000
$ws.Range("Q61").Value = "$/caja 16 kilos"
$ws.Range("R61").Value = "Provincia de Los Andes"
$ws.Range("S61").Value = 1375
$ws.Range("T61").Value = 16
